# Aggiunto indice Gulpease nuove versioni dei documenti
# Adds two new document-version rows ("1.1.0" and "2.0.0") to each of the
# five per-document Gulpease-index mini tables on the "Documentazione"
# sheet, updates a handful of existing "Verbali" scores in the overview
# table, extends the related charts' source ranges, and restores the
# sheet's last-used selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Documentazione")

# ---------------------------------------------------------------------
# 1) Overview table (rows 3-10): a few score corrections in column D
# ---------------------------------------------------------------------
$ws.Range("D4").Value = 71
$ws.Range("D9").Value = 59
$ws.Range("D10").Value = 61

# ---------------------------------------------------------------------
# 2) Helper: append the "1.1.0" / "2.0.0" rows right below an existing
#    mini table, copying the formatting of the table's last data row.
# ---------------------------------------------------------------------
function Add-VersionRows {
    param(
        [int]$LastRow,
        [int]$C110,
        [int]$C200
    )

    $src = $ws.Range("B$LastRow`:E$LastRow")
    $row1 = $LastRow + 1
    $row2 = $LastRow + 2

    $src.Copy() | Out-Null
    $ws.Range("B$row1`:E$row1").PasteSpecial(-4122) | Out-Null
    $src.Copy() | Out-Null
    $ws.Range("B$row2`:E$row2").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    $ws.Range("B$row1").Value = "1.1.0"
    $ws.Range("C$row1").Value = $C110
    $ws.Range("D$row1").Value = 40
    $ws.Range("E$row1").Value = 80

    $ws.Range("B$row2").Value = "2.0.0"
    $ws.Range("C$row2").Value = $C200
    $ws.Range("D$row2").Value = 40
    $ws.Range("E$row2").Value = 80
}

# Analisi dei requisiti : rows 60-63 -> add 64,65
Add-VersionRows -LastRow 63 -C110 60 -C200 65
# Norme di progetto     : rows 82-85 -> add 86,87
Add-VersionRows -LastRow 85 -C110 72 -C200 71
# Piano di progetto     : rows 105-107 -> add 108,109
Add-VersionRows -LastRow 107 -C110 70 -C200 70
# Piano di qualifica    : rows 124-126 -> add 127,128
Add-VersionRows -LastRow 126 -C110 65 -C200 68
# Glossario             : rows 167-169 -> add 170,171
Add-VersionRows -LastRow 169 -C110 60 -C200 59

# ---------------------------------------------------------------------
# 3) Extend the matching charts' source ranges so they include the two
#    new rows (categories + the three score series for each table).
# ---------------------------------------------------------------------
function Update-ChartRange {
    param(
        [int]$ChartIndex,
        [int]$HeaderRow,
        [int]$OldLastRow,
        [int]$NewLastRow
    )

    $chart = $ws.ChartObjects($ChartIndex).Chart
    $cols = @("C", "D", "E")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $seriesIndex = $i + 1
        $formula = "=SERIES(Documentazione!`$$col`$$HeaderRow,Documentazione!`$B`$$($HeaderRow + 1):`$B`$$NewLastRow,Documentazione!`$$col`$$($HeaderRow + 1):`$$col`$$NewLastRow,$seriesIndex)"
        $chart.SeriesCollection($seriesIndex).Formula = $formula
    }
}

Update-ChartRange -ChartIndex 2 -HeaderRow 59  -OldLastRow 63  -NewLastRow 65
Update-ChartRange -ChartIndex 3 -HeaderRow 81  -OldLastRow 85  -NewLastRow 87
Update-ChartRange -ChartIndex 4 -HeaderRow 104 -OldLastRow 107 -NewLastRow 109
Update-ChartRange -ChartIndex 5 -HeaderRow 123 -OldLastRow 126 -NewLastRow 128
Update-ChartRange -ChartIndex 7 -HeaderRow 166 -OldLastRow 169 -NewLastRow 171

# ---------------------------------------------------------------------
# 4) Restore the sheet view (scroll position / selection) left by the
#    author after the edit.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E49").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 107
